$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(99, 8).Value = 438.55554
$ws.Cells.Item(99, 9).Value = 388.91666
$ws.Cells.Item(99, 10).Value = 537.8333
$ws.Cells.Item(99, 11).Value = 1166.74998
$ws.Cells.Item(99, 12).Value = 1613.4999
$ws.Cells.Item(99, 13).Value = 331.2500199999999
$ws.Cells.Item(99, 14).Value = -4609.4999
$ws.Cells.Item(106, 8).Value = 1419.9
$ws.Cells.Item(106, 9).Value = 1244.3334
$ws.Cells.Item(106, 11).Value = 1244.3334
$ws.Cells.Item(106, 13).Value = -613.3334
$ws.Cells.Item(111, 8).Value = 1164
$ws.Cells.Item(111, 9).Value = 1164
$ws.Cells.Item(111, 11).Value = 3492
$ws.Cells.Item(111, 13).Value = -425
$ws.Cells.Item(116, 8).Value = 6142.4546
$ws.Cells.Item(116, 9).Value = 6248.857
$ws.Cells.Item(116, 10).Value = 5956.25
$ws.Cells.Item(116, 11).Value = 6248.857
$ws.Cells.Item(116, 12).Value = 5956.25
$ws.Cells.Item(116, 13).Value = -2806.857
$ws.Cells.Item(116, 14).Value = -12840.25
$ws.Cells.Item(132, 8).Value = 436339.34
$ws.Cells.Item(132, 9).Value = 1514.5238
$ws.Cells.Item(132, 11).Value = 4543.5714
$ws.Cells.Item(132, 13).Value = -2013.5714
$ws.Cells.Item(138, 8).Value = 2406.175
$ws.Cells.Item(138, 9).Value = 3653.1428
$ws.Cells.Item(138, 10).Value = 2141.6667
$ws.Cells.Item(138, 11).Value = 10959.4284
$ws.Cells.Item(138, 12).Value = 6425.000100000001
$ws.Cells.Item(138, 13).Value = -5819.428400000001
$ws.Cells.Item(138, 14).Value = -16705.0001
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 2431.3096
$ws.Cells.Item(32, 9).Value = 2052.975
$ws.Cells.Item(32, 11).Value = 2052.975
$ws.Cells.Item(32, 13).Value = -1765.975
$ws.Cells.Item(61, 8).Value = 2588.4736
$ws.Cells.Item(61, 9).Value = 2555.8572
$ws.Cells.Item(61, 11).Value = 2555.8572
$ws.Cells.Item(61, 13).Value = -2343.8572
$ws.Cells.Item(102, 8).Value = 2092.3333
$ws.Cells.Item(102, 9).Value = 2053.875
$ws.Cells.Item(102, 11).Value = 2053.875
$ws.Cells.Item(102, 13).Value = -431.875
$ws.Cells.Item(136, 8).Value = 2588.4736
$ws.Cells.Item(136, 9).Value = 2555.8572
$ws.Cells.Item(136, 11).Value = 7667.571599999999
$ws.Cells.Item(136, 13).Value = -5117.571599999999
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(96, 8).Value = 11048.714
$ws.Cells.Item(96, 9).Value = 11048.714
$ws.Cells.Item(96, 11).Value = 11048.714
$ws.Cells.Item(96, 13).Value = -8302.714
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 2484
$ws.Cells.Item(16, 9).Value = 2484
$ws.Cells.Item(16, 11).Value = 2484
$ws.Cells.Item(16, 13).Value = -2197
$ws.Cells.Item(22, 8).Value = 555
$ws.Cells.Item(22, 9).Value = 520
$ws.Cells.Item(22, 11).Value = 520
$ws.Cells.Item(22, 13).Value = -170
$ws.Cells.Item(58, 8).Value = 2579.3572
$ws.Cells.Item(58, 9).Value = 2857.6365
$ws.Cells.Item(58, 10).Value = 1559
$ws.Cells.Item(58, 11).Value = 2857.6365
$ws.Cells.Item(58, 12).Value = 1559
$ws.Cells.Item(58, 13).Value = -2654.6365
$ws.Cells.Item(58, 14).Value = -1965
$ws.Cells.Item(99, 8).Value = 12
$ws.Cells.Item(99, 9).Value = 12
$ws.Cells.Item(99, 11).Value = 12
$ws.Cells.Item(99, 13).Value = 1486
$ws.Cells.Item(107, 8).Value = 1719.7778
$ws.Cells.Item(107, 9).Value = 1719.7778
$ws.Cells.Item(107, 11).Value = 1719.7778
$ws.Cells.Item(107, 13).Value = 200.2221999999999
$ws.Cells.Item(113, 8).Value = 2484
$ws.Cells.Item(113, 9).Value = 2484
$ws.Cells.Item(113, 11).Value = 2484
$ws.Cells.Item(113, 13).Value = -314
$ws.Cells.Item(126, 8).Value = 12
$ws.Cells.Item(126, 9).Value = 12
$ws.Cells.Item(126, 11).Value = 36
$ws.Cells.Item(126, 13).Value = 2434
$ws.Cells.Item(136, 8).Value = 2579.3572
$ws.Cells.Item(136, 9).Value = 2857.6365
$ws.Cells.Item(136, 10).Value = 1559
$ws.Cells.Item(136, 11).Value = 8572.9095
$ws.Cells.Item(136, 12).Value = 4677
$ws.Cells.Item(136, 13).Value = -6022.9095
$ws.Cells.Item(136, 14).Value = -9777
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(38, 8).Value = 175.14285
$ws.Cells.Item(38, 9).Value = 46.333332
$ws.Cells.Item(38, 10).Value = 271.75
$ws.Cells.Item(38, 11).Value = 138.999996
$ws.Cells.Item(38, 12).Value = 815.25
$ws.Cells.Item(38, 13).Value = 208.000004
$ws.Cells.Item(38, 14).Value = -1509.25
$ws.Cells.Item(46, 8).Value = 17040.666
$ws.Cells.Item(46, 10).Value = 20428.8
$ws.Cells.Item(46, 12).Value = 61286.39999999999
$ws.Cells.Item(46, 14).Value = -61468.39999999999
$ws.Cells.Item(80, 8).Value = 10546.77
$ws.Cells.Item(80, 10).Value = 5811
$ws.Cells.Item(80, 12).Value = 17433
$ws.Cells.Item(80, 14).Value = -19305
$ws.Cells.Item(83, 8).Value = 10546.77
$ws.Cells.Item(83, 10).Value = 5811
$ws.Cells.Item(83, 12).Value = 52299
$ws.Cells.Item(83, 14).Value = -61659
$ws.Cells.Item(97, 8).Value = 911.4
$ws.Cells.Item(97, 9).Value = 1598.75
$ws.Cells.Item(97, 10).Value = 453.16666
$ws.Cells.Item(97, 11).Value = 4796.25
$ws.Cells.Item(97, 12).Value = 1359.49998
$ws.Cells.Item(97, 13).Value = -4300.25
$ws.Cells.Item(97, 14).Value = -2351.49998
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 3809.7778
$ws.Cells.Item(102, 9).Value = 3400
$ws.Cells.Item(102, 10).Value = 3926.8572
$ws.Cells.Item(102, 11).Value = 3400
$ws.Cells.Item(102, 12).Value = 3926.8572
$ws.Cells.Item(102, 13).Value = -1778
$ws.Cells.Item(102, 14).Value = -7170.8572
$ws.Cells.Item(113, 8).Value = 4458.5
$ws.Cells.Item(113, 9).Value = 3009.5
$ws.Cells.Item(113, 11).Value = 3009.5
$ws.Cells.Item(113, 13).Value = -839.5
$ws.Cells.Item(132, 8).Value = 2277.3
$ws.Cells.Item(132, 9).Value = 2277.3
$ws.Cells.Item(132, 11).Value = 6831.900000000001
$ws.Cells.Item(132, 13).Value = -4301.900000000001
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(61, 8).Value = 2700
$ws.Cells.Item(61, 9).Value = 2400
$ws.Cells.Item(61, 10).Value = 3000
$ws.Cells.Item(61, 11).Value = 2400
$ws.Cells.Item(61, 12).Value = 3000
$ws.Cells.Item(61, 13).Value = -2198
$ws.Cells.Item(61, 14).Value = -3404
$ws.Cells.Item(87, 8).Value = 25000
$ws.Cells.Item(87, 9).Value = 25000
$ws.Cells.Item(87, 11).Value = 25000
$ws.Cells.Item(87, 13).Value = -23877
$ws.Cells.Item(90, 8).Value = 25000
$ws.Cells.Item(90, 9).Value = 25000
$ws.Cells.Item(90, 11).Value = 75000
$ws.Cells.Item(90, 13).Value = -69384
$ws.Cells.Item(113, 8).Value = 2700
$ws.Cells.Item(113, 9).Value = 2400
$ws.Cells.Item(113, 10).Value = 3000
$ws.Cells.Item(113, 11).Value = 2400
$ws.Cells.Item(113, 12).Value = 3000
$ws.Cells.Item(113, 13).Value = -230
$ws.Cells.Item(113, 14).Value = -7340
$ws.Cells.Item(122, 8).Value = 3968.4443
$ws.Cells.Item(122, 9).Value = 4386.25
$ws.Cells.Item(122, 10).Value = 3360.7273
$ws.Cells.Item(122, 11).Value = 13158.75
$ws.Cells.Item(122, 12).Value = 10082.1819
$ws.Cells.Item(122, 13).Value = -10708.75
$ws.Cells.Item(122, 14).Value = -14982.1819
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(46, 8).Value = 27000
$ws.Cells.Item(46, 10).Value = 27000
$ws.Cells.Item(46, 12).Value = 27000
$ws.Cells.Item(46, 14).Value = -27462
$ws.Cells.Item(56, 8).Value = 27500
$ws.Cells.Item(56, 10).Value = 27500
$ws.Cells.Item(56, 12).Value = 27500
$ws.Cells.Item(56, 14).Value = -28928
$ws.Cells.Item(99, 8).Value = 34998
$ws.Cells.Item(99, 9).Value = 34998
$ws.Cells.Item(99, 11).Value = 34998
$ws.Cells.Item(99, 13).Value = -32003
$ws.Cells.Item(134, 8).Value = 27000
$ws.Cells.Item(134, 10).Value = 27000
$ws.Cells.Item(134, 12).Value = 81000
$ws.Cells.Item(134, 14).Value = -86070
